$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Helper: find the (1-based) index of the first paragraph whose
# trimmed text equals $text, searching within the supplied range of
# paragraph indices.
# ------------------------------------------------------------------
function Find-ParagraphIndex($doc, [string]$text, [int]$fromIndex, [int]$toIndex) {
    for ($i = $fromIndex; $i -le $toIndex; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.Trim() -eq $text) {
            return $i
        }
    }
    return -1
}

# ------------------------------------------------------------------
# 1. Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document. We clone the
#    plain-style (no pStyle) paragraph right after it as a template,
#    then overwrite its text, so the new paragraph picks up the
#    document's normal (unstyled) body formatting instead of
#    inheriting Heading1.
# ------------------------------------------------------------------
$titleIndex = Find-ParagraphIndex $d "Play Candyways Bonanza 2 Megaways Free: Slot Game Review" 1 $d.Paragraphs.Count
$titlePara = $d.Paragraphs.Item($titleIndex)

$templatePara = $d.Paragraphs.Item($titleIndex + 2)
$templatePara.Range.Copy()

$insertPoint = $d.Range($titlePara.Range.End, $titlePara.Range.End)
$insertPoint.Paste()

$metaPara = $d.Paragraphs.Item($titleIndex + 1)
$metaText = "Meta description: Get the lowdown on Candyways Bonanza 2 Megaways slot game with our review. Play for free and experience its variety of symbols, winning combos, and bonuses."
$metaPara.Range.Text = $metaText

# Bold just the "Meta description" label (leave the rest, including the
# colon, in regular formatting).
$labelStart = $metaPara.Range.Start
$labelRange = $d.Range($labelStart, $labelStart + 16)
$labelRange.Font.Bold = 1

# ------------------------------------------------------------------
# 2. At the bottom of the document, remove the duplicated bold title
#    paragraph and turn the remaining italic paragraph into the
#    "Feature Image Prompt" paragraph.
# ------------------------------------------------------------------
$dupTitleIndex = Find-ParagraphIndex $d "Play Candyways Bonanza 2 Megaways Free: Slot Game Review" ($titleIndex + 1) $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs.Item($dupTitleIndex)
$dupTitlePara.Range.Delete()

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)

$featureText = "Feature Image Prompt: Design a cartoon-style feature image for Candyways Bonanza 2 Megaways that features a happy Maya warrior with glasses. The happy warrior should be holding a giant candy wand with the game's logo at the top. The warrior should be standing in front of a colorful and vibrant background filled with candy and sweets. The image should have a fun, playful vibe that will appeal to players who enjoy colorful and visually appealing online slot games."
$lastRange.Text = $featureText

Write-Host "Meta description paragraph inserted at index" ($titleIndex + 1)
Write-Host "Duplicate title paragraph removed (was index" $dupTitleIndex ")"
Write-Host "Final paragraph count:" $d.Paragraphs.Count
